$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 47109.2
$ws.Range("I32").Value = 142766.33
$ws.Range("J32").Value = 6113.2856
$ws.Range("K32").Value = 142766.33
$ws.Range("L32").Value = 6113.2856
$ws.Range("M32").Value = -142440.33
$ws.Range("N32").Value = -6765.2856
$ws.Range("H62").Value = 31537
$ws.Range("I62").Value = 4450.5
$ws.Range("K62").Value = 4450.5
$ws.Range("M62").Value = -3826.5
$ws.Range("H65").Value = 31537
$ws.Range("I65").Value = 4450.5
$ws.Range("K65").Value = 22252.5
$ws.Range("M65").Value = -19132.5
$ws.Range("H111").Value = 45706.5
$ws.Range("I111").Value = 2261
$ws.Range("J111").Value = 74670.164
$ws.Range("K111").Value = 6783
$ws.Range("L111").Value = 224010.492
$ws.Range("M111").Value = -3716
$ws.Range("N111").Value = -230144.492
$ws.Range("H113").Value = 76927180
$ws.Range("J113").Value = 4257
$ws.Range("L113").Value = 4257
$ws.Range("N113").Value = -10765
$ws.Range("H121").Value = 1587.8889
$ws.Range("J121").Value = 1587.8889
$ws.Range("L121").Value = 4763.6667
$ws.Range("N121").Value = -8257.6667
$ws.Range("H132").Value = 1993.4736
$ws.Range("I132").Value = 1514.3715
$ws.Range("K132").Value = 4543.1145
$ws.Range("M132").Value = -2013.1145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2175
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1123
$ws.Range("H61").Value = 3311.9565
$ws.Range("I61").Value = 1920.625
$ws.Range("K61").Value = 1920.625
$ws.Range("M61").Value = -1708.625
$ws.Range("H63").Value = 70593080
$ws.Range("I63").Value = 166668580
$ws.Range("K63").Value = 166668580
$ws.Range("M63").Value = -166667894
$ws.Range("H66").Value = 70593080
$ws.Range("I66").Value = 166668580
$ws.Range("K66").Value = 833342900
$ws.Range("M66").Value = -833339468
$ws.Range("H88").Value = 15153246
$ws.Range("J88").Value = 2449.8333
$ws.Range("L88").Value = 2449.8333
$ws.Range("N88").Value = -3261.8333
$ws.Range("H91").Value = 15153246
$ws.Range("J91").Value = 2449.8333
$ws.Range("L91").Value = 2449.8333
$ws.Range("N91").Value = -5257.8333
$ws.Range("H110").Value = 1344.76
$ws.Range("I110").Value = 1078.2273
$ws.Range("J110").Value = 3299.3333
$ws.Range("K110").Value = 1078.2273
$ws.Range("L110").Value = 3299.3333
$ws.Range("M110").Value = 966.7727
$ws.Range("N110").Value = -7389.3333
$ws.Range("H122").Value = 3493.077
$ws.Range("I122").Value = 2082
$ws.Range("J122").Value = 4375
$ws.Range("K122").Value = 6246
$ws.Range("L122").Value = 13125
$ws.Range("M122").Value = -3796
$ws.Range("N122").Value = -18025
$ws.Range("H132").Value = 6986.615
$ws.Range("I132").Value = 7168.3
$ws.Range("K132").Value = 21504.9
$ws.Range("M132").Value = -18974.9
$ws.Range("H136").Value = 3311.9565
$ws.Range("I136").Value = 1920.625
$ws.Range("K136").Value = 5761.875
$ws.Range("M136").Value = -3211.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 2560
$ws.Range("I29").Value = 2264
$ws.Range("J29").Value = 3300
$ws.Range("K29").Value = 2264
$ws.Range("L29").Value = 3300
$ws.Range("M29").Value = -1975
$ws.Range("N29").Value = -3878
$ws.Range("H86").Value = 1569.8889
$ws.Range("I86").Value = 1353.5588
$ws.Range("J86").Value = 5247.5
$ws.Range("K86").Value = 1353.5588
$ws.Range("L86").Value = 5247.5
$ws.Range("M86").Value = -230.5588
$ws.Range("N86").Value = -7493.5
$ws.Range("H89").Value = 1569.8889
$ws.Range("I89").Value = 1353.5588
$ws.Range("J89").Value = 5247.5
$ws.Range("K89").Value = 6767.794
$ws.Range("L89").Value = 26237.5
$ws.Range("M89").Value = -1151.794
$ws.Range("N89").Value = -37469.5
$ws.Range("H126").Value = 119946.25
$ws.Range("J126").Value = 119946.25
$ws.Range("L126").Value = 119946.25
$ws.Range("N126").Value = -129826.25
$ws.Range("H134").Value = 2802.8333
$ws.Range("I134").Value = 2242.0732
$ws.Range("J134").Value = 6087.2856
$ws.Range("K134").Value = 6726.219599999999
$ws.Range("L134").Value = 18261.8568
$ws.Range("M134").Value = -4191.219599999999
$ws.Range("N134").Value = -23331.8568

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2954.6052
$ws.Range("I31").Value = 1156.8334
$ws.Range("K31").Value = 1156.8334
$ws.Range("M31").Value = -861.8334
$ws.Range("H34").Value = 2954.6052
$ws.Range("I34").Value = 1156.8334
$ws.Range("K34").Value = 1156.8334
$ws.Range("M34").Value = -954.8334
$ws.Range("H99").Value = 11300
$ws.Range("J99").Value = 11100
$ws.Range("L99").Value = 11100
$ws.Range("N99").Value = -14096
$ws.Range("H122").Value = 2483
$ws.Range("J122").Value = 2132.6667
$ws.Range("L122").Value = 6398.000100000001
$ws.Range("N122").Value = -11298.0001
$ws.Range("H126").Value = 11300
$ws.Range("J126").Value = 11100
$ws.Range("L126").Value = 33300
$ws.Range("N126").Value = -38240

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 31883.75
$ws.Range("I58").Value = 33997
$ws.Range("K58").Value = 33997
$ws.Range("M58").Value = -33720
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H132").Value = 3087.8235
$ws.Range("I132").Value = 3027.5806
$ws.Range("J132").Value = 3710.3333
$ws.Range("K132").Value = 9082.7418
$ws.Range("L132").Value = 11130.9999
$ws.Range("M132").Value = -6552.7418
$ws.Range("N132").Value = -16190.9999
$ws.Range("H136").Value = 26393.072
$ws.Range("J136").Value = 26393.072
$ws.Range("L136").Value = 79179.216
$ws.Range("N136").Value = -84279.216

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3172.8462
$ws.Range("I7").Value = 2856.2666
$ws.Range("J7").Value = 3604.5454
$ws.Range("K7").Value = 2856.2666
$ws.Range("L7").Value = 3604.5454
$ws.Range("M7").Value = -2744.2666
$ws.Range("N7").Value = -3828.5454
$ws.Range("H126").Value = 3172.8462
$ws.Range("I126").Value = 2856.2666
$ws.Range("J126").Value = 3604.5454
$ws.Range("K126").Value = 8568.799800000001
$ws.Range("L126").Value = 10813.6362
$ws.Range("M126").Value = -6098.799800000001
$ws.Range("N126").Value = -15753.6362
$ws.Range("H136").Value = 6443.2583
$ws.Range("I136").Value = 8177.143
$ws.Range("K136").Value = 24531.429
$ws.Range("M136").Value = -21981.429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 24729
$ws.Range("J41").Value = 24729
$ws.Range("L41").Value = 24729
$ws.Range("N41").Value = -25509
$ws.Range("H58").Value = 19995.5
$ws.Range("I58").Value = 19995.5
$ws.Range("K58").Value = 19995.5
$ws.Range("M58").Value = -19687.5
$ws.Range("H107").Value = 882.4706
$ws.Range("I107").Value = 834.1818
$ws.Range("J107").Value = 971
$ws.Range("K107").Value = 2502.5454
$ws.Range("L107").Value = 2913
$ws.Range("M107").Value = -582.5454
$ws.Range("N107").Value = -6753
$ws.Range("H132").Value = 5856.7
$ws.Range("I132").Value = 5669.0835
$ws.Range("K132").Value = 17007.2505
$ws.Range("M132").Value = -14477.2505
$ws.Range("H136").Value = 6687.9165
$ws.Range("I136").Value = 8621.6
$ws.Range("J136").Value = 3465.111
$ws.Range("K136").Value = 25864.8
$ws.Range("L136").Value = 10395.333
$ws.Range("M136").Value = -23314.8
$ws.Range("N136").Value = -15495.333
